# Plantilla_Cargue_Productos.xlsx
# - Remove the "N° Lote" and "Fecha Vencimiento" columns (originally D and E).
# - Insert a new "Precio Costo" column between "Precio Venta" and "Cantidad".
# - Update the selected cell to F4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns D (N° Lote) and E (Fecha Vencimiento)
$ws.Range("D1:E1").EntireColumn.Delete()

# After the deletion the header row is:
#   A=Nombre  B=Presentacion  C=Laboratorio  D=Precio Venta  E=Cantidad
# Insert a new column before the (now) "Cantidad" column for "Precio Costo"
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("E1").Value = "Precio Costo"
$ws.Columns.Item(5).ColumnWidth = 11.5

# Match the saved selection from the authored workbook
$ws.Range("F4").Select()
